$d = $word.ActiveDocument

# Locate the paragraph that reads "Second pass revision of the Psalter."
# (the first bullet in the TODO list) and the paragraph immediately
# following it ("Move the Canticles ..."), which carries the _GoBack
# bookmark at its end.
$firstIndex = 0
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "Second pass revision of the Psalter.*") {
        $firstIndex = $i
        break
    }
}

if ($firstIndex -eq 0) {
    throw "Could not find the 'Second pass revision of the Psalter.' paragraph"
}

$secondIndex = $firstIndex + 1

$p1 = $d.Paragraphs.Item($firstIndex)
$r1 = $p1.Range
$p2 = $d.Paragraphs.Item($secondIndex)
$r2 = $p2.Range

# 1. Remove the text of the first paragraph ("Second pass revision of the
#    Psalter. ") but keep its paragraph mark, so the paragraph (and the
#    formatting / identity it carries) survives.
$textRange = $d.Range($r1.Start, $r1.End - 1)
$textRange.Delete()

# 2. Re-resolve the (now empty) first paragraph and the second paragraph,
#    whose text we are about to fold into the first paragraph.
$p1b = $d.Paragraphs.Item($firstIndex)
$r1b = $p1b.Range
$p2b = $d.Paragraphs.Item($secondIndex)
$r2b = $p2b.Range

$secondTextOnly = $d.Range($r2b.Start, $r2b.End - 1)
$insertionStart = $r1b.Start
$insertionPoint = $d.Range($insertionStart, $insertionStart)
$insertionPoint.InsertAfter($secondTextOnly.Text)

# 3. Delete the now-redundant second paragraph (its text has been copied
#    into the first paragraph; this removes its own paragraph mark too).
$p2c = $d.Paragraphs.Item($secondIndex)
$r2c = $p2c.Range
$r2c.Delete()

# 4. Recreate the _GoBack bookmark at the start of the merged paragraph
#    (it used to sit at the end of the second paragraph; in the merged
#    paragraph it now sits right before the run).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
$bmRange = $d.Range($insertionStart, $insertionStart)
$d.Bookmarks.Add("_GoBack", $bmRange)
